# Add a new "Serviced by " column (O) to the Card23 sheet, and fix the
# trailing-space typo in the existing "Correction " header (N1 -> "Correction").
# Also backfill the "nan" placeholder text into the (previously blank) N2:N12
# cells, matching the rest of the row's placeholder values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")

# Fix header typo: "Correction " (trailing space) -> "Correction"
$ws.Range("N1").Value = "Correction"

# New header for the added column
$ws.Range("O1").Value = "Serviced by "

# Copy the header formatting (bold, centered, bordered) from N1 onto the
# newly added O1 header cell so it matches the rest of the header row.
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Backfill "nan" placeholder text that the original data rows use for blank
# values into the previously-empty N2:N12 cells.
$ws.Range("N2:N12").Value = "nan"
